$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Accuracy")

$updates = @{
    "V2"  = 0.2118068965517241
    "V3"  = 0.1805458229957766
    "V7"  = 0.2757242757242757
    "V11" = 0.216893039049236
    "V12" = 0.2738805263656158
    "V14" = 0.2939434534301629
    "V15" = 0.2832591683289857
    "V16" = 0.2672127950068266
    "V17" = 0.2971181376820615
    "V22" = 0.2408601378302409
    "V27" = 0.2075536661425302
    "V28" = 0.2561708439026856
    "V29" = 0.2733596891596734
    "V30" = 0.2227304056875329
    "V32" = 0.2399536512231055
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
